$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "Finding Text (FIND)" formula in column N so that it correctly
# extracts the region name (e.g. "West" or "North") regardless of its
# length, instead of always grabbing a fixed 4 characters (which truncated
# "North" to "Nort").
$ws.Range("N4").Formula = "=MID(K4,4,FIND("" "",K4)-4)"
$ws.Range("N5:N38").Formula = "=MID(K5,4,FIND("" "",K5)-4)"

# Reflect the cell the author ended up with selected after editing the
# formula.
$ws.Activate()
$ws.Range("N4").Select()
